# Finalized Experiments with Participant Generation
# Renames each task-order worksheet and refreshes the generated file-id values.

$wb = $excel.ActiveWorkbook

# --- Rename worksheets (new participant-generation run ids) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16502911469084547"
$wb.Worksheets.Item(2).Name = "NB_TO-165029114866179"
$wb.Worksheets.Item(3).Name = "RS_TO-16502911486627953"
$wb.Worksheets.Item(4).Name = "TOL_TO-1650291148720528"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16502911487985556"

# --- Sheet 1: GNG ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16502911468745914.csv"
$ws1.Range("B3").Value = "GNG_stims-16502911468912878.csv"
$ws1.Range("B4").Value = "go_stims-1650291146892322.csv"
$ws1.Range("B5").Value = "GNG_stims-16502911469074512.csv"

# --- Sheet 2: NB ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-1650291147791087.csv"
$ws2.Range("B3").Value = "TB-1650291148260164.csv"
$ws2.Range("B4").Value = "ZB-match_6-16502911469745142.csv"
$ws2.Range("B5").Value = "TB-16502911479695494.csv"
$ws2.Range("B6").Value = "ZB-match_7-165029114719781.csv"
$ws2.Range("B7").Value = "OB-16502911478653367.csv"
$ws2.Range("B8").Value = "ZB-match_1-16502911472707343.csv"
$ws2.Range("B9").Value = "TB-1650291148638005.csv"
$ws2.Range("B10").Value = "OB-16502911478812156.csv"

# --- Sheet 3: RS ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes open"
$ws3.Range("B3").Value = "eyes closed"

# --- Sheet 4: TOL ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16502911486773612.csv"
$ws4.Range("B3").Value = "ZM_stims-16502911486647913.csv"
$ws4.Range("B4").Value = "MM_stims-16502911487035546.csv"
$ws4.Range("B5").Value = "ZM_stims-16502911486783614.csv"
$ws4.Range("B6").Value = "MM_stims-16502911487195306.csv"
$ws4.Range("B7").Value = "ZM_stims-16502911487045622.csv"

# --- Sheet 5: vSAT ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16502911487225702.csv"
$ws5.Range("B3").Value = "SAT_stims-1650291148735407.csv"
$ws5.Range("B4").Value = "vSAT_stims-16502911487840965.csv"
$ws5.Range("B5").Value = "vSAT_stims-1650291148752206.csv"
